$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-7 : TICKER, DATE, BUY/SELL, PRICE, VOLUME, NET_EFFECT_TO_CASH,
#                 TOTAL_SHARES_HOLDING, TICKER_TOTAL_VALUE, AVERAGE_PRICE, REALIZED_PROFIT
$data = @(
    @("AAPL", "08/14/20", "BUY",  100,    3000,  -300000,    3000,  300000,     100,   ""),
    @("AAPL", "08/14/20", "BUY",  200,    3000,  -600000,    6000,  900000,     150,   ""),
    @("AAPL", "08/14/20", "SELL", 500,    6000,  3000000,    0,     0,          150,   ""),
    @("CODX", "08/14/20", "BUY",  300,    100,   -30000,     100,   30000,      300,   ""),
    @("CODX", "08/14/20", "BUY",  23.32,  12999, -303136.68, 13099, 333136.68,  25.43, ""),
    @("CODX", "08/14/20", "SELL", 300,    3000,  900000,     10099, 256846.68,  25.43, "")
)

# Make sure the DATE column (B) is treated as plain text, not auto-converted
# into a date serial number by Excel's cell-value inference.
$ws.Range("B2:B7").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($value in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $value
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
